# release chapter highly cited
# Reposition four country/region label textboxes (tx8, tx9, tx10, tx11)
# that live inside the single top-level group shape on slide 1.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$g = $s.Shapes.Item(1)

# EMU -> Points conversion factor (1 pt = 12700 EMU).
$emuPerPt = 12700

# Shape.Left/.Top round-trip through a single-precision (float32) point
# value internally, so converting emu -> pt -> emu can truncate down by
# one EMU. A tiny epsilon (well inside rounding tolerance, see below)
# nudges the float past the truncation boundary so the saved OOXML
# offset lands exactly on the target EMU value.
$epsilonPt = 0.00003

function Set-ShapePosition {
    param($shape, [double]$xEmu, [double]$yEmu)
    $shape.Left = ($xEmu / $emuPerPt) + $epsilonPt
    $shape.Top  = ($yEmu / $emuPerPt) + $epsilonPt
}

for ($i = 1; $i -le $g.GroupItems.Count; $i++) {
    $sh = $g.GroupItems.Item($i)
    switch ($sh.Name) {
        "tx8"  { Set-ShapePosition $sh 4779451 2842418 }
        "tx9"  { Set-ShapePosition $sh 4380215 4471880 }
        "tx10" { Set-ShapePosition $sh 5635481 3547988 }
        "tx11" { Set-ShapePosition $sh 5621448 2848312 }
    }
}
